$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data contained in row 18 and row 19 (all columns), which effectively
# exchanges the two occurrence records while keeping unchanged shared fields as-is.

$row18 = 18
$row19 = 19

# Capture current (pre-edit) values from row 18
$A18 = $ws.Cells.Item($row18, 1).Value2
$B18 = $ws.Cells.Item($row18, 2).Value2
$D18 = $ws.Cells.Item($row18, 4).Value2
$E18 = $ws.Cells.Item($row18, 5).Value2
$F18 = $ws.Cells.Item($row18, 6).Value2
$G18 = $ws.Cells.Item($row18, 7).Value2
$H18 = $ws.Cells.Item($row18, 8).Value2
$Q18 = $ws.Cells.Item($row18, 17).Value2
$R18 = $ws.Cells.Item($row18, 18).Value2
$Z18 = $ws.Cells.Item($row18, 26).Value2
$AB18 = $ws.Cells.Item($row18, 28).Value2
$AC18 = $ws.Cells.Item($row18, 29).Value2

# Capture current (pre-edit) values from row 19
$A19 = $ws.Cells.Item($row19, 1).Value2
$B19 = $ws.Cells.Item($row19, 2).Value2
$D19 = $ws.Cells.Item($row19, 4).Value2
$E19 = $ws.Cells.Item($row19, 5).Value2
$F19 = $ws.Cells.Item($row19, 6).Value2
$G19 = $ws.Cells.Item($row19, 7).Value2
$H19 = $ws.Cells.Item($row19, 8).Value2
$Q19 = $ws.Cells.Item($row19, 17).Value2
$R19 = $ws.Cells.Item($row19, 18).Value2
$Z19 = $ws.Cells.Item($row19, 26).Value2
$AB19 = $ws.Cells.Item($row19, 28).Value2

# Write row 19's original values into row 18
$ws.Cells.Item($row18, 1).Value2 = $A19
$ws.Cells.Item($row18, 2).Value2 = $B19
$ws.Cells.Item($row18, 4).Value2 = $D19
$ws.Cells.Item($row18, 5).Value2 = $E19
$ws.Cells.Item($row18, 6).Value2 = $F19
$ws.Cells.Item($row18, 7).Value2 = $G19
$ws.Cells.Item($row18, 8).Value2 = $H19
$ws.Cells.Item($row18, 17).Value2 = $Q19
$ws.Cells.Item($row18, 18).Value2 = $R19
$ws.Cells.Item($row18, 26).Value2 = $Z19
$ws.Cells.Item($row18, 28).Value2 = $AB19
$ws.Cells.Item($row18, 29).ClearContents()

# Write row 18's original values into row 19
$ws.Cells.Item($row19, 1).Value2 = $A18
$ws.Cells.Item($row19, 2).Value2 = $B18
$ws.Cells.Item($row19, 4).Value2 = $D18
$ws.Cells.Item($row19, 5).Value2 = $E18
$ws.Cells.Item($row19, 6).Value2 = $F18
$ws.Cells.Item($row19, 7).Value2 = $G18
$ws.Cells.Item($row19, 8).Value2 = $H18
$ws.Cells.Item($row19, 17).Value2 = $Q18
$ws.Cells.Item($row19, 18).Value2 = $R18
$ws.Cells.Item($row19, 26).Value2 = $Z18
$ws.Cells.Item($row19, 28).Value2 = $AB18
$ws.Cells.Item($row19, 29).Value2 = $AC18
